# "ora mi mette tutte le parole chiave affianco"
# Replace the D column ("Nome" -> "ParoleChiave") with a longer, richer list
# of keywords, drop the stray A2 value, and extend the sheet with new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook default font size 12 -> 11 (Normal cell style)
$wb.Styles.Item("Normal").Font.Size = 11

# Header rename: "Nome" -> "ParoleChiave"
$ws.Range("D1").Value = "ParoleChiave"

# Drop the leftover process-number in A2 (only A1 keeps a value now)
$ws.Range("A2").Value = $null

# Updated/expanded keyword column (D) for the existing rows
$ws.Range("D2").Value = "packaging sostenibile"
$ws.Range("D3").Value = "imballaggio sostenibile"
$ws.Range("D4").Value = "packaging ecologico"
$ws.Range("D5").Value = "imballaggio ecologico"
$ws.Range("D6").Value = "packaging biodegradabile"
$ws.Range("D7").Value = "packaging compostabile"
$ws.Range("D8").Value = "packaging riciclabile"
$ws.Range("D9").Value = "imballaggio riciclabile"

# New rows 10-24: sequential index in B, empty search text in C, new keyword in D
$newRows = @(
    @{ Row = 10; Index = "9";  Keyword = "carta kraft" },
    @{ Row = 11; Index = "10"; Keyword = "carta riciclata" },
    @{ Row = 12; Index = "11"; Keyword = "cellulosa di bambù" },
    @{ Row = 13; Index = "12"; Keyword = "fibra di bambù" },
    @{ Row = 14; Index = "13"; Keyword = "materiale riciclato" },
    @{ Row = 15; Index = "14"; Keyword = "materiale ecologico" },
    @{ Row = 16; Index = "15"; Keyword = "materiale sostenibile" },
    @{ Row = 17; Index = "16"; Keyword = "bambù" },
    @{ Row = 18; Index = "17"; Keyword = "campione gratuito" },
    @{ Row = 19; Index = "18"; Keyword = "spedizione campioni" },
    @{ Row = 20; Index = "19"; Keyword = "campione personalizzato" },
    @{ Row = 21; Index = "20"; Keyword = "packaging personalizzato" },
    @{ Row = 22; Index = "21"; Keyword = "imballaggio personalizzato" },
    @{ Row = 23; Index = "22"; Keyword = "scatola ecologica" },
    @{ Row = 24; Index = "23"; Keyword = "scatola sostenibile" }
)

foreach ($r in $newRows) {
    $cellB = $ws.Cells.Item($r.Row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.Index
    $cellB.Style = "Normal"

    $ws.Cells.Item($r.Row, 4).Value = $r.Keyword
}

